$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").Style = "Normal"
$ws.Range("N2").Value = "2019-09-30 00:00:00"

$ws.Range("O2").Value = 848185811.38
$ws.Range("P2").Value = 353.0199665662
$ws.Range("Q2").Value = 26428232505.7
$ws.Range("R2").Value = 10999.5871545974
$ws.Range("S2").Value = 1409076132.92
$ws.Range("T2").Value = 586.465845878
$ws.Range("U2").Value = -757211343.98
$ws.Range("V2").Value = -315.1558535275
$ws.Range("W2").Value = 40201341.11
$ws.Range("X2").Value = 16.7320366648
$ws.Range("Y2").Value = 177312484.28
$ws.Range("Z2").Value = 73.79850786519999
$ws.Range("AA2").Value = -331240144.05
$ws.Range("AB2").Value = -137.8641130387
$ws.Range("AC2").Value = -240265676.65
$ws.Range("AD2").Value = -142.1303191631
